{"js": "// Peer-review update: replace the OR / pval / moderator values in the\n// single results table with the corrected figures.\n//\n// The table has a fixed, known shape (9 rows x 3 cols: moderators | OR |\n// pval), so each changed value is addressed directly by its (row, col)\n// (0-indexed, per Office.js) rather than via a document-wide search. We\n// assert the existing cell text matches what we expect, then replace just\n// that text run (via a search scoped to the cell, so only the matched\n// range's text/run is touched and the rest of the run formatting is kept\n// intact) with the corrected value.\nconst replacements = [\n  // row, col, oldText, newText\n  [1, 1, \"5530691178.68 (0.01-3.79160731696921e+21)\", \"0.23 (0-1252664.13)\"],\n  [1, 2, \"0.10667188\", \"0.8510603\"],\n  [2, 1, \"0.61 (0.36-1.04)\", \"0.99 (0.73-1.34)\"],\n  [2, 2, \"0.07066908\", \"0.9393696\"],\n  [3, 1, \"1.03 (0.97-1.09)\", \"1.01 (0.94-1.09)\"],\n  [3, 2, \"0.30161683\", \"0.7649440\"],\n  [4, 1, \"1 (0.96-1.05)\", \"1 (0.97-1.03)\"],\n  [4, 2, \"0.87534944\", \"0.8985404\"],\n  [5, 1, \"0.87 (0.57-1.33)\", \"1.21 (0.91-1.61)\"],\n  [5, 2, \"0.50877734\", \"0.1885343\"],\n  [6, 0, \"Behavioural support only\", \"EMA study type - Interventional\"],\n  [6, 1, \"17.74 (0.48-659.03)\", \"0.31 (0.04-2.67)\"],\n  [6, 2, \"0.11888934\", \"0.2838084\"],\n  [7, 0, \"Combined support\", \"Study duration days\"],\n  [7, 1, \"228.08 (0.69-75144.33)\", \"0.98 (0.88-1.09)\"],\n  [7, 2, \"0.06641412\", \"0.6848491\"],\n  [8, 2, \"0.61547937\", \"0.6154795\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nfor (const [row, col, oldText, newText] of replacements) {\n  const cell = table.getCell(row, col);\n  cell.body.load(\"text\");\n  await context.sync();\n\n  if (cell.body.text !== oldText) {\n    throw new Error(\n      `Unexpected text in cell (${row},${col}): ${JSON.stringify(cell.body.text)} expected ${JSON.stringify(oldText)}`\n    );\n  }\n\n  const results = cell.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Peer-review update: replace the OR / pval / moderator values in the\n# single results table with the corrected figures.\n#\n# The table has a fixed, known shape (9 rows x 3 cols: moderators | OR |\n# pval), so we address each changed value directly by (row, col) via the\n# Word table object model rather than a document-wide Find/Replace. Each\n# cell's Range is trimmed to drop the trailing end-of-cell mark before its\n# .Text is reassigned, and we assert the existing text matches what we\n# expect before overwriting it, as a guard against addressing the wrong\n# cell.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nfunction Set-CellText($row, $col, $old, $new) {\n    $cell = $tbl.Cell($row, $col)\n    $r = $cell.Range\n    $r.End = $r.End - 1\n    if ($r.Text -ne $old) {\n        throw \"Unexpected text in cell ($row,$col): [$($r.Text)] expected [$old]\"\n    }\n    $r.Text = $new\n}\n\nSet-CellText 2 2 \"5530691178.68 (0.01-3.79160731696921e+21)\" \"0.23 (0-1252664.13)\"\nSet-CellText 2 3 \"0.10667188\" \"0.8510603\"\nSet-CellText 3 2 \"0.61 (0.36-1.04)\" \"0.99 (0.73-1.34)\"\nSet-CellText 3 3 \"0.07066908\" \"0.9393696\"\nSet-CellText 4 2 \"1.03 (0.97-1.09)\" \"1.01 (0.94-1.09)\"\nSet-CellText 4 3 \"0.30161683\" \"0.7649440\"\nSet-CellText 5 2 \"1 (0.96-1.05)\" \"1 (0.97-1.03)\"\nSet-CellText 5 3 \"0.87534944\" \"0.8985404\"\nSet-CellText 6 2 \"0.87 (0.57-1.33)\" \"1.21 (0.91-1.61)\"\nSet-CellText 6 3 \"0.50877734\" \"0.1885343\"\nSet-CellText 7 1 \"Behavioural support only\" \"EMA study type - Interventional\"\nSet-CellText 7 2 \"17.74 (0.48-659.03)\" \"0.31 (0.04-2.67)\"\nSet-CellText 7 3 \"0.11888934\" \"0.2838084\"\nSet-CellText 8 1 \"Combined support\" \"Study duration days\"\nSet-CellText 8 2 \"228.08 (0.69-75144.33)\" \"0.98 (0.88-1.09)\"\nSet-CellText 8 3 \"0.06641412\" \"0.6848491\"\nSet-CellText 9 3 \"0.61547937\" \"0.6154795\"\n\nWrite-Output \"done\"\n"}
